$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value looks like a plain number (e.g. "253.89") must be
# force-formatted as Text first, otherwise Excel auto-converts the assigned
# string into a numeric value (losing the original text-cell semantics/
# formatting such as preserved trailing zeros).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"

$ws.Range("D2").Value = "35.123.60"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.904.78"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D5").Value = "253.89"
$ws.Range("E5").Value = "  +3.27%  "
$ws.Range("D6").Value = "0.701"
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("D8").Value = "41.54"
$ws.Range("E8").Value = "  +2.72%  "
$ws.Range("E9").Value = "  +4.17%  "
$ws.Range("D10").Value = "52.43"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").Value = "0.0753"
$ws.Range("D12").Value = "0.0979"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "13.21"
$ws.Range("E13").Value = "  +5.06%  "
$ws.Range("D14").Value = "2.181.96"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("E15").Value = "  +4.65%  "
$ws.Range("E16").Value = "  +5.52%  "
$ws.Range("D17").Value = "1.906.94"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "35.118.72"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "73.83"
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("E20").Value = "  +3.16%  "
$ws.Range("D21").Value = "243.13"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").Value = "13.02"
$ws.Range("E22").Value = "  +4.00%  "
$ws.Range("D23").Value = "5.06"
$ws.Range("E23").Value = "  +6.21%  "
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("E25").Value = "  +5.50%  "
$ws.Range("D26").Value = "2.32"
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("D27").Value = "167.68"
$ws.Range("D28").Value = "8.59"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "4.128.16"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("E32").Value = "  +7.40%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "2.02"
$ws.Range("E33").Value = "  +8.05%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "4.34"
$ws.Range("E34").Value = "  +4.94%  "
$ws.Range("E35").Value = "  +7.82%  "
$ws.Range("E36").Value = "  +3.45%  "
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").Value = "0.854"
$ws.Range("E38").Value = "  -5.52%  "
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("D40").Value = "103.08"
$ws.Range("E40").Value = "  +15.22%  "
$ws.Range("D41").Value = "17.24"
$ws.Range("E41").Value = "  +7.09%  "
$ws.Range("E42").Value = "  +4.13%  "
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("D44").Value = "0.0650"
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("D45").Value = "2.42"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").Value = "1.302.41"
$ws.Range("E46").Value = "  -2.95%  "
$ws.Range("D47").Value = "12.70"
$ws.Range("E47").Value = "  +4.67%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("E51").Value = "  +7.24%  "
